# Implement login page ui elements
# Append 3 new data rows (106-108) to each of the 4 worksheets in the
# FE2025 database workbook, matching the existing row layout/format.

$wb = $excel.ActiveWorkbook

# Data block per worksheet (1-based worksheet index): each entry is
# @(A_date, B, C, D, E, F, G, H, I)
$data = @{
    1 = @(
        @(45892.49715277777, "0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x00,0xF8", "0xf", 380, 759863127514710900000000.0, 248, 15),
        @(45893.495,         "0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x00,0xF4", "0xf", 380, 759863127514710900000000.0, 244, 15),
        @(45894.4950925926,  "0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x00,0xF4", "0xf", 380, 759863127514710900000000.0, 244, 15)
    )
    2 = @(
        @(45892.49715277777, "0x01,0x90", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x01,0x04", "0xe", 400, 568432987514711000000000.0, 260, 14),
        @(45893.495,         "0x01,0x90", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x01,0x00", "0xe", 400, 568432987514711000000000.0, 256, 14),
        @(45894.4950925926,  "0x01,0x90", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x01,0x00", "0xe", 400, 568432987514711000000000.0, 256, 14)
    )
    3 = @(
        @(45892.49715277777, "0x00,0x6e", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x5A", "0x3", 110, 568631262647114000000000.0, 90, 3),
        @(45893.495,         "0x00,0x6e", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x5A", "0x3", 110, 568631262647114000000000.0, 90, 3),
        @(45894.4950925926,  "0x00,0x6e", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x5A", "0x3", 110, 568631262647114000000000.0, 90, 3)
    )
    4 = @(
        @(45892.49715277777, "0x00,0x6e", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x59", "0x3", 110, 985046333984776000000000.0, 89, 3),
        @(45893.495,         "0x00,0x6e", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x58", "0x3", 110, 985046333984776000000000.0, 88, 3),
        @(45894.4950925926,  "0x00,0x6e", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x58", "0x3", 110, 985046333984776000000000.0, 88, 3)
    )
}

foreach ($sheetIndex in 1..4) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $rows = $data[$sheetIndex]
    $startRow = 106

    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $startRow + $i
        $row = $rows[$i]

        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 6).Value = $row[5]
        $ws.Cells.Item($r, 7).Value = $row[6]
        $ws.Cells.Item($r, 8).Value = $row[7]
        $ws.Cells.Item($r, 9).Value = $row[8]
    }
}
